$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Correspond Handoff Datetime (D3) and
# Correspond Handback DateTime (G3) to reflect the new report generation times.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-11 03:28:39"
$wsZhCn.Range("G3").Value = "2016-01-11 03:29:26"

# "de-de" sheet: same columns updated.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-11 03:28:53"
$wsDeDe.Range("G3").Value = "2016-01-11 03:29:52"
